$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the data value: B15 changes from 5.7 to 6
$ws.Range("B15").Value = 6

# Move the active cell selection to G15 (matches the user's last selection)
$ws.Range("G15").Select()
